# Updated cryptos list (refreshed price/volume snapshot), plus a couple of
# rank swaps: rows 22/23 swap Litecoin/Dai, and row 51 swaps Arweave -> ONDO.
# Leading "'" forces text entry so numeric-looking strings (e.g. "582.02",
# "68.107.01") stay stored as text, matching the source inlineStr cells
# instead of being auto-converted to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.107.01"
$ws.Range("E2").Value = "'  +0.25%  "
$ws.Range("D3").Value = "'3.250.17"
$ws.Range("E3").Value = "'  -0.16%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'582.02"
$ws.Range("E5").Value = "'  +0.50%  "
$ws.Range("D6").Value = "'184.89"
$ws.Range("E6").Value = "'  +1.08%  "
$ws.Range("E8").Value = "'  +0.65%  "
$ws.Range("E9").Value = "'  -2.84%  "
$ws.Range("E10").Value = "'  -1.15%  "
$ws.Range("D11").Value = "'0.417"
$ws.Range("E11").Value = "'  +0.40%  "
$ws.Range("D12").Value = "'3.810.83"
$ws.Range("E12").Value = "'  -0.20%  "
$ws.Range("E13").Value = "'  +0.05%  "
$ws.Range("D14").Value = "'28.00"
$ws.Range("E14").Value = "'  -2.67%  "
$ws.Range("D15").Value = "'68.090.90"
$ws.Range("E15").Value = "'  +0.23%  "
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("E16").Value = "'  -0.35%  "
$ws.Range("D17").Value = "'3.278.10"
$ws.Range("E17").Value = "'  +0.80%  "
$ws.Range("D18").Value = "'5.81"
$ws.Range("E18").Value = "'  -0.47%  "
$ws.Range("D19").Value = "'13.50"
$ws.Range("E19").Value = "'  -0.28%  "
$ws.Range("D20").Value = "'395.65"
$ws.Range("E20").Value = "'  +4.37%  "
$ws.Range("D21").Value = "'7.60"
$ws.Range("E21").Value = "'  -0.65%  "
$ws.Range("B22").Value = "'Dai"
$ws.Range("C22").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "'  +0.10%  "
$ws.Range("B23").Value = "'Litecoin"
$ws.Range("C23").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").Value = "'71.38"
$ws.Range("E23").Value = "'  +0.14%  "
$ws.Range("E24").Value = "'  +0.70%  "
$ws.Range("E25").Value = "'  -0.81%  "
$ws.Range("E26").Value = "'  +3.08%  "
$ws.Range("D27").Value = "'9.67"
$ws.Range("E27").Value = "'  -4.01%  "
$ws.Range("E29").Value = "'  -0.21%  "
$ws.Range("D30").Value = "'5.63"
$ws.Range("E30").Value = "'  -0.75%  "
$ws.Range("D31").Value = "'22.81"
$ws.Range("E31").Value = "'  -0.18%  "
$ws.Range("D32").Value = "'7.05"
$ws.Range("E32").Value = "'  +0.31%  "
$ws.Range("E33").Value = "'  -0.07%  "
$ws.Range("E34").Value = "'  +0.11%  "
$ws.Range("E35").Value = "'  -5.18%  "
$ws.Range("D36").Value = "'162.10"
$ws.Range("E36").Value = "'  -0.32%  "
$ws.Range("D37").Value = "'1.91"
$ws.Range("E37").Value = "'  +2.41%  "
$ws.Range("E38").Value = "'  -3.26%  "
$ws.Range("D39").Value = "'4.62"
$ws.Range("E39").Value = "'  +0.49%  "
$ws.Range("D40").Value = "'26.48"
$ws.Range("E40").Value = "'  -0.08%  "
$ws.Range("E41").Value = "'  -1.11%  "
$ws.Range("E42").Value = "'  -4.58%  "
$ws.Range("D43").Value = "'41.02"
$ws.Range("E43").Value = "'  -0.28%  "
$ws.Range("D44").Value = "'0.0685"
$ws.Range("E44").Value = "'  -0.18%  "
$ws.Range("D45").Value = "'25.14"
$ws.Range("E45").Value = "'  -1.66%  "
$ws.Range("D46").Value = "'2.611.05"
$ws.Range("E46").Value = "'  -0.97%  "
$ws.Range("D47").Value = "'338.06"
$ws.Range("E47").Value = "'  -2.63%  "
$ws.Range("D48").Value = "'0.0279"
$ws.Range("E48").Value = "'  -2.56%  "
$ws.Range("D49").Value = "'6.35"
$ws.Range("E49").Value = "'  +3.07%  "
$ws.Range("E50").Value = "'  -0.73%  "
$ws.Range("B51").Value = "'ONDO"
$ws.Range("C51").Value = "'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Value = "'0.982"
$ws.Range("E51").Value = "'  -1.13%  "
